# Update 'want to go' counts (column F) per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 8138
$ws.Range("F4").Value = 1914
$ws.Range("F5").Value = 6507
$ws.Range("F7").Value = 2064
$ws.Range("F8").Value = 567
$ws.Range("F9").Value = 44
$ws.Range("F11").Value = 48
$ws.Range("F15").Value = 8503
$ws.Range("F16").Value = 159
$ws.Range("F20").Value = 1805
$ws.Range("F25").Value = 30
$ws.Range("F30").Value = 2072
$ws.Range("F31").Value = 846
$ws.Range("F32").Value = 469
$ws.Range("F35").Value = 175
$ws.Range("F36").Value = 145
$ws.Range("F37").Value = 3
$ws.Range("F38").Value = 25
$ws.Range("F40").Value = 3962

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 390
$ws.Range("F21").Value = 36

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 390
$ws.Range("F5").Value = 8138
$ws.Range("F8").Value = 1914
$ws.Range("F9").Value = 6507
$ws.Range("F11").Value = 2064
$ws.Range("F14").Value = 567
$ws.Range("F15").Value = 44
$ws.Range("F18").Value = 48
$ws.Range("F23").Value = 8503
$ws.Range("F24").Value = 159
$ws.Range("F27").Value = 1805
$ws.Range("F32").Value = 2072
$ws.Range("F33").Value = 846
$ws.Range("F35").Value = 469
$ws.Range("F43").Value = 3962
$ws.Range("F44").Value = 36
